# Apply the edits described by the diff:
#  - sharedStrings: patient01@gmail.com -> alice.brown@gmail.com (G2 display text)
#  - H2 contact number: 91234567 -> 99991234
#  - selection moves from G2 to H2
#  - workbook theme colours updated to the new "Office" (Aptos) palette
#  - workbook theme fonts updated to Aptos Display / Aptos Narrow

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value changes -----------------------------------------------
$ws.Range("G2").Value = "alice.brown@gmail.com"
$ws.Range("H2").Value = 99991234

# --- Selection moves to H2 ---------------------------------------------
$ws.Range("H2").Select() | Out-Null

# --- Theme colour scheme -------------------------------------------------
# Colors(1)=dk1 Colors(2)=lt1 Colors(3)=dk2 Colors(4)=lt2
# Colors(5..10)=accent1..accent6 Colors(11)=hlink Colors(12)=folHlink
$tcs = $wb.Theme.ThemeColorScheme
$tcs.Colors(3).RGB  = 4270094    # dk2      0E2841
$tcs.Colors(4).RGB  = 15263976   # lt2      E8E8E8
$tcs.Colors(5).RGB  = 8544277    # accent1  156082
$tcs.Colors(6).RGB  = 3305961    # accent2  E97132
$tcs.Colors(7).RGB  = 2386713    # accent3  196B24
$tcs.Colors(8).RGB  = 13999631   # accent4  0F9ED5
$tcs.Colors(9).RGB  = 9644960    # accent5  A02B93
$tcs.Colors(10).RGB = 3057486    # accent6  4EA72E
$tcs.Colors(11).RGB = 8812614    # hlink    467886
$tcs.Colors(12).RGB = 8216726    # folHlink 96607D

# --- Theme font scheme ----------------------------------------------------
$tfs = $wb.Theme.ThemeFontScheme
$tfs.MajorFont.Name = "Aptos Display"
$tfs.MinorFont.Name = "Aptos Narrow"
